# Rubrica Andamento - registrando configurações de ambiente de deploy
# Updates the "Porcentagem" (F) column for the items that moved to
# "Concluido" status, and leaves the view positioned on C19 (matching the
# last selection made by the author while reviewing the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escopo")

# Subcompetencias that got finished -> bring their progress to 100%.
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("F18").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("F25").Value = 1

# F26 (TOTAL CONCLUIDO) is a live AVERAGE(F2:F25) formula, so it recalculates
# automatically from 60% to ~86% once the cells above are updated.

# Leave the sheet scrolled/selected the way the author left it.
$ws.Activate() | Out-Null
$ws.Range("C19").Select() | Out-Null
